$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value = "ARWXAB"
$ws.Range("B14").Value = "Chip contact para cartucho de tóner"
$ws.Range("C14").Value = "MFP136NW 108A 108W 103A 103W 131A 133PN 136A 138P 138PN 138PNW"
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 70000
$ws.Range("F14").Value = 7
$ws.Range("G14").Value = 2
$ws.Range("H14").Formula = "=(E14-D14)*G14"
$ws.Range("I14").Formula = "=D14*F14"
$ws.Range("J14").Value = 0
